# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.246.80"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.853.83"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.78"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4607"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3704"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8858"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.05"
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07808"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.384"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.776.31"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.41"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008929"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.74"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "27.257.69"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.105"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "2.070.41"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.925"
$ws.Range("E25").Value = "  +4.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.42"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.056"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.84"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.058"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08803"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.091"
$ws.Range("E32").Value = "  +4.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7677"
$ws.Range("E33").Value = "  +4.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.170"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.500"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.733"
$ws.Range("E36").Value = "  +10.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.082"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01948"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05251"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.066"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5105"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1628"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.385"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4789"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.34"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.34"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.641"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06217"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.71"
$ws.Range("E51").Value = "  +1.43%  "
